$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 373.9091
$ws.Range("I2").Value = 331.4
$ws.Range("K2").Value = 331.4
$ws.Range("M2").Value = -218.4
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H32").Value = 6368.8335
$ws.Range("I32").Value = 4880.1665
$ws.Range("J32").Value = 7857.5
$ws.Range("K32").Value = 4880.1665
$ws.Range("L32").Value = 7857.5
$ws.Range("M32").Value = -4554.1665
$ws.Range("N32").Value = -8509.5
$ws.Range("H80").Value = 1428774.9
$ws.Range("I80").Value = 7609223.5
$ws.Range("J80").Value = 2517.4614
$ws.Range("K80").Value = 22827670.5
$ws.Range("L80").Value = 7552.3842
$ws.Range("M80").Value = -22826672.5
$ws.Range("N80").Value = -9548.3842
$ws.Range("H83").Value = 1428774.9
$ws.Range("I83").Value = 7609223.5
$ws.Range("J83").Value = 2517.4614
$ws.Range("K83").Value = 68483011.5
$ws.Range("L83").Value = 22657.1526
$ws.Range("M83").Value = -68478019.5
$ws.Range("N83").Value = -32641.1526
$ws.Range("H88").Value = 3666.6667
$ws.Range("J88").Value = 3666.6667
$ws.Range("L88").Value = 3666.6667
$ws.Range("N88").Value = -4478.6667
$ws.Range("H91").Value = 3666.6667
$ws.Range("J91").Value = 3666.6667
$ws.Range("L91").Value = 3666.6667
$ws.Range("N91").Value = -6474.6667
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H111").Value = 13781.223
$ws.Range("J111").Value = 20934
$ws.Range("L111").Value = 62802
$ws.Range("N111").Value = -68936
$ws.Range("H129").Value = 1436.9474
$ws.Range("I129").Value = 1200.8572
$ws.Range("K129").Value = 3602.5716
$ws.Range("M129").Value = 1397.4284
$ws.Range("H131").Value = 6168.4614
$ws.Range("I131").Value = 1979
$ws.Range("J131").Value = 20133.334
$ws.Range("K131").Value = 5937
$ws.Range("L131").Value = 60400.00199999999
$ws.Range("M131").Value = -897
$ws.Range("N131").Value = -70480.00199999999
$ws.Range("H132").Value = 3652.5322
$ws.Range("I132").Value = 3479.5833
$ws.Range("J132").Value = 4245.5
$ws.Range("K132").Value = 10438.7499
$ws.Range("L132").Value = 12736.5
$ws.Range("M132").Value = -7908.749899999999
$ws.Range("N132").Value = -17796.5
$ws.Range("H133").Value = 131778
$ws.Range("J133").Value = 131778
$ws.Range("L133").Value = 131778
$ws.Range("N133").Value = -141898

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 38724.555
$ws.Range("I31").Value = 15896.6
$ws.Range("J31").Value = 67259.5
$ws.Range("K31").Value = 15896.6
$ws.Range("L31").Value = 67259.5
$ws.Range("M31").Value = -15602.6
$ws.Range("N31").Value = -67847.5
$ws.Range("H32").Value = 176339.11
$ws.Range("I32").Value = 188181.17
$ws.Range("K32").Value = 188181.17
$ws.Range("M32").Value = -187894.17
$ws.Range("H45").Value = 1616
$ws.Range("I45").Value = 1462.8
$ws.Range("K45").Value = 1462.8
$ws.Range("M45").Value = -1085.8
$ws.Range("H61").Value = 28573472
$ws.Range("I61").Value = 30305016
$ws.Range("J61").Value = 2999.5
$ws.Range("K61").Value = 30305016
$ws.Range("L61").Value = 2999.5
$ws.Range("M61").Value = -30304804
$ws.Range("N61").Value = -3423.5
$ws.Range("H132").Value = 2618
$ws.Range("I132").Value = 2068.4614
$ws.Range("K132").Value = 6205.3842
$ws.Range("M132").Value = -3675.3842
$ws.Range("H136").Value = 28573472
$ws.Range("I136").Value = 30305016
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 90915048
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -90912498
$ws.Range("N136").Value = -14098.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5055.5
$ws.Range("I94").Value = 5083.25
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 5083.25
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -4632.25
$ws.Range("N94").Value = -5902
$ws.Range("H97").Value = 2809.3333
$ws.Range("I97").Value = 464
$ws.Range("J97").Value = 7500
$ws.Range("K97").Value = 464
$ws.Range("L97").Value = 7500
$ws.Range("M97").Value = 527
$ws.Range("N97").Value = -9482
$ws.Range("H105").Value = 3028.15
$ws.Range("I105").Value = 3544.5386
$ws.Range("K105").Value = 3544.5386
$ws.Range("M105").Value = -1797.5386
$ws.Range("H134").Value = 889.0862
$ws.Range("I134").Value = 797.3261
$ws.Range("K134").Value = 2391.9783
$ws.Range("M134").Value = 143.0217000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 50332.668
$ws.Range("J97").Value = 46999
$ws.Range("L97").Value = 46999
$ws.Range("N97").Value = -48981
$ws.Range("H132").Value = 1054358.6
$ws.Range("I132").Value = 909920.4
$ws.Range("J132").Value = 1252961.2
$ws.Range("K132").Value = 2729761.2
$ws.Range("L132").Value = 3758883.6
$ws.Range("M132").Value = -2727231.2
$ws.Range("N132").Value = -3763943.6
$ws.Range("H134").Value = 1852.2667
$ws.Range("I134").Value = 982.0833
$ws.Range("J134").Value = 5333
$ws.Range("K134").Value = 2946.2499
$ws.Range("L134").Value = 15999
$ws.Range("M134").Value = -411.2498999999998
$ws.Range("N134").Value = -21069

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2999.5
$ws.Range("J116").Value = 2999.5
$ws.Range("L116").Value = 8998.5
$ws.Range("N116").Value = -15882.5
$ws.Range("H129").Value = 2263.6
$ws.Range("J129").Value = 2583.4
$ws.Range("L129").Value = 7750.200000000001
$ws.Range("N129").Value = -17750.2
$ws.Range("H131").Value = 4871.9443
$ws.Range("I131").Value = 1449.1666
$ws.Range("J131").Value = 6583.3335
$ws.Range("K131").Value = 4347.4998
$ws.Range("L131").Value = 19750.0005
$ws.Range("M131").Value = 692.5002000000004
$ws.Range("N131").Value = -29830.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14079.728
$ws.Range("I70").Value = 13875.777
$ws.Range("K70").Value = 13875.777
$ws.Range("M70").Value = -13605.777
$ws.Range("H73").Value = 14079.728
$ws.Range("I73").Value = 13875.777
$ws.Range("K73").Value = 13875.777
$ws.Range("M73").Value = -12939.777
$ws.Range("H128").Value = 49999.5
$ws.Range("J128").Value = 49999.5
$ws.Range("L128").Value = 49999.5
$ws.Range("N128").Value = -59959.5
$ws.Range("H132").Value = 422112
$ws.Range("I132").Value = 560644.4
$ws.Range("J132").Value = 6514.8335
$ws.Range("K132").Value = 1681933.2
$ws.Range("L132").Value = 19544.5005
$ws.Range("M132").Value = -1679403.2
$ws.Range("N132").Value = -24604.5005
$ws.Range("H141").Value = 90999.5
$ws.Range("J141").Value = 90999.5
$ws.Range("L141").Value = 90999.5
$ws.Range("N141").Value = -101359.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 893.55554
$ws.Range("J22").Value = 941.8333
$ws.Range("L22").Value = 941.8333
$ws.Range("N22").Value = -1531.8333
$ws.Range("H27").Value = 893.55554
$ws.Range("J27").Value = 941.8333
$ws.Range("L27").Value = 941.8333
$ws.Range("N27").Value = -1155.8333
$ws.Range("H70").Value = 24461
$ws.Range("I70").Value = 8777
$ws.Range("K70").Value = 8777
$ws.Range("M70").Value = -8507
$ws.Range("H73").Value = 24461
$ws.Range("I73").Value = 8777
$ws.Range("K73").Value = 8777
$ws.Range("M73").Value = -7841
$ws.Range("H136").Value = 1837.1154
$ws.Range("I136").Value = 1632.4791
$ws.Range("K136").Value = 4897.4373
$ws.Range("M136").Value = -2347.4373

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16677583
$ws.Range("I81").Value = 7500.5
$ws.Range("J81").Value = 20011600
$ws.Range("K81").Value = 15001
$ws.Range("L81").Value = 40023200
$ws.Range("M81").Value = -13940
$ws.Range("N81").Value = -40025322
$ws.Range("H84").Value = 16677583
$ws.Range("I84").Value = 7500.5
$ws.Range("J84").Value = 20011600
$ws.Range("K84").Value = 75005
$ws.Range("L84").Value = 200116000
$ws.Range("M84").Value = -69701
$ws.Range("N84").Value = -200126608
$ws.Range("H136").Value = 3302.6
$ws.Range("I136").Value = 2485.12
$ws.Range("K136").Value = 7455.36
$ws.Range("M136").Value = -4905.36
$ws.Range("H140").Value = 98999.5
$ws.Range("J140").Value = 98999.5
$ws.Range("L140").Value = 98999.5
$ws.Range("N140").Value = -109359.5
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360
